# Update the existing performance table with new values and append
# additional rows, keeping all numeric-looking values stored as text
# (matching the workbook's existing convention of text-typed number cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows starting at row 2: Name, LODA value, LOF value
$data = @(
    @("Annthyroid", "0.0593", "0.1967"),
    @("Arrhythmia", "0.1753", "0.3798"),
    @("Breastw", "0.6431", "0.3453"),
    @("Glass", "0.0411", "0.1092"),
    @("Ionosphere", "0.7711", "0.8635"),
    @("Letter", "0.113", "0.2714"),
    @("Lympho", "0.2946", "0.8012"),
    @("Mammography", "0.1886", "0.1381"),
    @("Mnist", "0.107", "0.3401"),
    @("Musk", "0.138", "0.0836"),
    @("Optdigits", "0.0172", "0.0222"),
    @("Pendigits", "0.1309", "0.0282"),
    @("Pima", "0.5441", "0.4686"),
    @("Satellite", "0.2217", "0.3958"),
    @("SatImage-2", "0.5139", "0.0422"),
    @("Shuttle", "0.4371", "0.123"),
    @("Speech", "0.0184", "0.0194"),
    @("Thyroid", "0.0151", "0.2832"),
    @("Vertebral", "0.0886", "0.0847"),
    @("Vowels", "0.0274", "0.4071"),
    @("Wbc", "0.4221", "0.5965"),
    @("Wine", "0.633", "0.3367")
)

# Force columns B and C to text format so numeric-looking strings are
# not reinterpreted as numbers (mirrors the source file's t="str" cells).
$lastRow = 1 + $data.Count
$ws.Range("B2:C$lastRow").NumberFormat = "@"

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Extend the "numbers stored as text" ignored-error marker over the full
# table (A1:C<lastRow>), matching how it previously covered A1:C3.
try {
    $fullRange = $ws.Range("A1:C$lastRow")
    $fullRange.Errors(3).Ignore = $true
} catch {
}
